$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "28.388.52"
Set-TextValue "E2" "  +0.07%  "
Set-TextValue "D3" "1.819.82"
Set-TextValue "E3" "  -0.55%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "315.06"
Set-TextValue "E5" "  -0.81%  "
Set-TextValue "E6" "  +0.01%  "
Set-TextValue "D7" "0.5104"
Set-TextValue "E7" "  -4.36%  "
Set-TextValue "D8" "0.3921"
Set-TextValue "E8" "  -3.33%  "
Set-TextValue "D9" "0.07805"
Set-TextValue "E9" "  +2.46%  "
Set-TextValue "D10" "41.74"
Set-TextValue "E10" "  -0.28%  "
Set-TextValue "D11" "1.108"
Set-TextValue "E11" "  +0.19%  "
Set-TextValue "D12" "20.91"
Set-TextValue "E12" "  +0.22%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "6.243"
Set-TextValue "E13" "  -1.78%  "
Set-TextValue "B14" "BinanceUSD"
Set-TextValue "C14" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D14" "1.001"
Set-TextValue "E14" "  +0.04%  "
Set-TextValue "D15" "7.474"
Set-TextValue "E15" "  -1.58%  "
Set-TextValue "D16" "1.815.68"
Set-TextValue "E16" "  -0.67%  "
Set-TextValue "D17" "0.00001137"
Set-TextValue "E17" "  +5.74%  "
Set-TextValue "D18" "92.45"
Set-TextValue "E18" "  +3.41%  "
Set-TextValue "D19" "0.06618"
Set-TextValue "E19" "  -0.05%  "
Set-TextValue "D20" "17.68"
Set-TextValue "E20" "  +0.30%  "
Set-TextValue "E21" "  -0.01%  "
Set-TextValue "D22" "6.079"
Set-TextValue "E22" "  -0.42%  "
Set-TextValue "D23" "28.427.44"
Set-TextValue "E23" "  +0.09%  "
Set-TextValue "D24" "11.24"
Set-TextValue "E24" "  +0.15%  "
Set-TextValue "D25" "2.242"
Set-TextValue "E25" "  +3.42%  "
Set-TextValue "D26" "21.13"
Set-TextValue "E26" "  +2.54%  "
Set-TextValue "D27" "2.023.41"
Set-TextValue "E27" "  -0.85%  "
Set-TextValue "D28" "155.09"
Set-TextValue "E28" "  -1.35%  "
Set-TextValue "D29" "2.395"
Set-TextValue "E29" "  -3.32%  "
Set-TextValue "D30" "125.26"
Set-TextValue "E30" "  +0.68%  "
Set-TextValue "E31" "  +0.54%  "
Set-TextValue "D32" "1.101"
Set-TextValue "E32" "  -1.99%  "
Set-TextValue "D33" "5.660"
Set-TextValue "E33" "  -0.61%  "
Set-TextValue "D34" "3.647"
Set-TextValue "E34" "  +0.32%  "
Set-TextValue "D35" "0.07046"
Set-TextValue "E35" "  -1.60%  "
Set-TextValue "D36" "0.2210"
Set-TextValue "E36" "  -2.27%  "
Set-TextValue "D37" "0.02320"
Set-TextValue "E37" "  -0.99%  "
Set-TextValue "D38" "5.173"
Set-TextValue "E38" "  -0.80%  "
Set-TextValue "D39" "8.751"
Set-TextValue "E39" "  -0.95%  "
Set-TextValue "D40" "0.6255"
Set-TextValue "E40" "  -0.27%  "
Set-TextValue "D41" "11.20"
Set-TextValue "E41" "  -1.06%  "
Set-TextValue "E42" "  -1.16%  "
Set-TextValue "D43" "0.9998"
Set-TextValue "E43" "  -0.07%  "
Set-TextValue "D44" "1.391"
Set-TextValue "E44" "  -0.66%  "
Set-TextValue "D45" "13.45"
Set-TextValue "E45" "  -0.82%  "
Set-TextValue "D46" "3.729"
Set-TextValue "E46" "  +0.66%  "
Set-TextValue "D47" "0.5875"
Set-TextValue "E47" "  +0.28%  "
Set-TextValue "D48" "124.18"
Set-TextValue "E48" "  -1.19%  "
Set-TextValue "D49" "1.973"
Set-TextValue "E49" "  -0.84%  "
Set-TextValue "E50" "  -0.91%  "
Set-TextValue "D51" "0.06897"
Set-TextValue "E51" "  +0.04%  "
